$d = $word.ActiveDocument

# 1. Change the title text
$d.Content.Find.Execute("Webinar 1 Write-up (Or put title here)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Replication Packages for Journals: For and Against", 2)

# 2. Change style of "Why Shouldn't Journals Require Replication Packages?" heading from Heading2 to Heading3
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    if ($text -like "Why Shouldn't *") {
        $p.Range.set_Style("Heading 3")
    }
}
